$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '62.849.23'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  -0.79%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '3.447.21'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  -1.12%  '
$ws.Range("E4").Value = '  +0.00%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '579.39'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -0.84%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '147.96'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +0.23%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  +0.32%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '7.92'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +2.98%  '
$ws.Range("E10").Value = '  -2.45%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.406'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +2.48%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '4.038.99'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -1.07%  '
$ws.Range("E13").Value = '  +2.63%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '28.24'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -4.88%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '3.451.79'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -1.00%  '
$ws.Range("E16").Value = '  -1.28%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '62.897.76'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -0.73%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '6.46'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +2.84%  '
$ws.Range("E19").Value = '  +1.47%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '9.12'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -2.78%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '388.96'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -0.61%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '0.561'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -0.36%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '74.83'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -0.04%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '3.593.46'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -3.06%  '
$ws.Range("E27").Value = '  -0.35%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '7.64'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -2.41%  '
$ws.Range("E29").Value = '  -0.07%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '8.03'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -3.03%  '
$ws.Range("E31").Value = '  -1.70%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("E33").Value = '  -5.83%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '23.29'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -2.08%  '
$ws.Range("E35").Value = '  +3.50%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '5.32'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -0.35%  '
$ws.Range("B37").Value = 'EnergySwap'
$ws.Range("C37").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '31.83'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -2.32%  '
$ws.Range("B38").Value = 'Aptos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '7.01'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -1.72%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '170.16'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -0.81%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '3.485.05'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -1.07%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.0783'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +2.01%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.791'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -2.00%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '42.66'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  +0.47%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '1.71'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -1.22%  '
$ws.Range("E45").Value = '  -3.31%  '
$ws.Range("E46").Value = '  -3.00%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '2.561.71'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -2.56%  '
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '6.89'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +1.96%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '2.26'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -1.57%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '22.64'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -4.31%  '
$ws.Range("E51").Value = '  +0.00%  '
